$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Poroto granado" at
# Vega Monumental Concepción. It belongs chronologically right after the
# existing row 6 (2021-05-26), so insert a new row at position 7 and push
# the rest of the data (old rows 7-23) down to rows 8-24.
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = "Poroto granado"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 26000
$ws.Range("M7").Value = 25500
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 1020
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
